# fix: reports fixed - fix getStudentDeadlineReport Error
#
# The header cell A1 currently reads "Name"; it must be changed to
# "Student Name" (the data/other headers are unchanged).
# Also restore the last active selection to G7, matching the saved file.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

$ws.Range("A1").Value = "Student Name"

$ws.Activate()
$ws.Range("G7").Select()
